# Update TAP playsheets/reports: DHMSM -> MHS GENESIS, per commit:
# "Changes to tap custom playsheets and reports to reflect changes in TAP
#  databases from DHMSM to MHS Genesis and for Disposition property changes"

$wb = $excel.ActiveWorkbook

function Repeat-Char($ch, $count) {
    $result = ""
    for ($i = 0; $i -lt $count; $i++) {
        $result = $result + $ch
    }
    return $result
}

$nbsp = [char]160
$rsquo = [char]8217

# ---------------------------------------------------------------------------
# "Additional Questions" sheet - Section 1 (Test Preparation) questions b-f
# ---------------------------------------------------------------------------
$wsAQ = $wb.Worksheets.Item("Additional Questions")

$wsAQ.Range("A6").Value = "b. What configuration activities would need to be performed to integrate this system" + $rsquo + "s test environment with the MHS GENESIS Test Data Center?"

$leadSpace = Repeat-Char $nbsp 52
$midSpace = Repeat-Char $nbsp 4
$wsAQ.Range("A7").Value = $leadSpace + "i." + $midSpace + " Has your team been coordinating with the MHS GENESIS Test & Evaluation group to plan for integration with the MHS GENESIS Test Data Center?"

$wsAQ.Range("A8").Value = "c. Does the program currently have funding available to develop a new interface with MHS GENESIS?"

$wsAQ.Range("A9").Value = "d. Does the program currently have contract scope to perform the work required to develop a new interface with MHS GENESIS?"

$wsAQ.Range("A10").Value = "e. What is the level of effort and/or timeline in which the program can be modernized to support the MHS GENESIS interface requirements, as specified in the Legacy Systems Modernization Tasker?"

$wsAQ.Range("A11").Value = "f. Are there any other risks that could affect this system" + $rsquo + "s ability to be modernized and configured to support MHS GENESIS T&E activities?"

# ---------------------------------------------------------------------------
# "Report Overview" sheet - report description blurb
# ---------------------------------------------------------------------------
$wsRO = $wb.Worksheets.Item("Report Overview")

$reportDetails = "This report details the characteristics and modernization activities required to transition the LPI system to the future-state environment. The report provides the following information and analysis regarding the LPI system:`n" + `
"1. System Overview - includes the system description along with any known Points of Contacts.`n" + `
"2. Additional Questions - Please provide/review answers to the questions from MHS GENESIS T&E or ICWG on this sheet. `n" + `
"3. Current-State Interface Definition - includes the systems current interfaces. Confirm the interfaces that the system currently supports, including the data objects exchanged, format, frequency and protocol; identify additional interface characteristics for the system, including the trigger events, ports, availability and response time`n" + `
"4. Future-State Interfaces for Development, Decommissioning, and Sustainment - based upon the capability defined in the MHS GENESIS Requirements Traceability Matrix, validate the interfaces that are being proposed and will need to be supported by the system; review the list of systems expected to be replaced by MHS GENESIS and validate that all data gaps are addressed"

$wsRO.Range("A8").Value = $reportDetails

# ---------------------------------------------------------------------------
# View-state: move the active tab/selection back to "Report Overview",
# leaving "Additional Questions" cursor parked on A9 (last-edited cell).
# ---------------------------------------------------------------------------
$wsAQ.Activate()
$wsAQ.Range("A9").Select() | Out-Null

$wsRO.Activate()
